$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "自己" (sheet1 / index 1): add column C = DEC2HEX(B) for rows 1-9.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C1").Formula = "=DEC2HEX(B1)"
$ws1.Range("C2:C9").Formula = "=DEC2HEX(B2)"

# ---------------------------------------------------------------------------
# Sheet "其他" (sheet2 / index 2): add column C = DEC2HEX(B) for row 1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C1").Formula = "=DEC2HEX(B1)"

# ---------------------------------------------------------------------------
# Sheet "黄宝辉" (sheet3 / index 3): add column C = DEC2HEX(B) for existing
# rows 1-3, plus three brand new rows (4-6) of data + formula.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C1").Formula = "=DEC2HEX(B1)"
$ws3.Range("C2:C3").Formula = "=DEC2HEX(B2)"

$ws3.Range("A4").Value = "7A114"
$ws3.Range("B4").Value = 21896
$ws3.Range("C4").Formula = "=DEC2HEX(B4)"

$ws3.Range("A5").Value = "7A104"
$ws3.Range("B5").Value = 21888
$ws3.Range("C5:C6").Formula = "=DEC2HEX(B5)"

$ws3.Range("A6").Value = "7A215"
$ws3.Range("B6").Value = 21911

$ws3.PageSetup.Orientation = 1

# Restore sheet3's own selection (new data pushed the cursor around).
$ws3.Range("C16").Select()

# ---------------------------------------------------------------------------
# Selections / active sheet: move the active tab from "其他" back to "自己",
# and leave each sheet's own cursor where the author left it.
# ---------------------------------------------------------------------------
$ws2.Range("C1").Select()

$ws1.Activate()
$ws1.Range("C1:C9").Select()
